$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: split "Compressive Estimation of Millimeter-Wave Channels" into
#    "Hybrid Digital and Analog Beamforming Design" / "for Large-Scale
#    Antenna Arrays" on two separate (but identically formatted) paragraphs.
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$found = $titleRange.Find.Execute("Compressive Estimation of Millimeter-Wave Channels")
if ($found) {
    $titleRange.Text = "Hybrid Digital and Analog Beamforming Design"
    # Splits the paragraph: a new (identically-formatted) paragraph is
    # created right after the current one, ready to receive the 2nd line.
    $titleRange.InsertParagraphAfter()

    $titlePara = $titleRange.Paragraphs.Item(1)
    $secondLinePara = $d.Paragraphs.Item($titlePara.Index + 1)
    $secondLinePara.Range.Text = "for Large-Scale Antenna Arrays"
}

# ---------------------------------------------------------------------------
# 2) Drop the stale <w:lastRenderedPageBreak/> markers that used to sit at
#    the very start of each section heading's first run. Re-keying the
#    first run's own text (blank it, then put the original characters back)
#    forces the run to be re-emitted without the cached page-break marker
#    while leaving its visible text untouched.
# ---------------------------------------------------------------------------
function Find-ParagraphByText($doc, $text) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $text) {
            return $p
        }
    }
    return $null
}

function Clear-LeadingPageBreak($doc, $paragraphText, $firstRunLength) {
    $p = Find-ParagraphByText $doc $paragraphText
    if ($p -eq $null) {
        return
    }
    $start = $p.Range.Start
    $runRange = $doc.Range($start, $start + $firstRunLength)
    $original = $runRange.Text
    # Re-use the SAME Range object for both writes: Word ranges track their
    # own content, so after shrinking it to "." it still points at exactly
    # that (now 1-char) span, and writing the original text back over it
    # restores the text while forcing the run (and its cached
    # lastRenderedPageBreak marker) to be re-emitted.
    $runRange.Text = "."
    $runRange.Text = $original
}

Clear-LeadingPageBreak $d "Abstract" 1
Clear-LeadingPageBreak $d "System Model" 12
Clear-LeadingPageBreak $d "Main Part" 4
Clear-LeadingPageBreak $d "Results and Discussion" 23
Clear-LeadingPageBreak $d "Conclusion" 10
